# Apply updated TPM values to the Col4a5-Cd93 LR-pairs sheet.
# Target cluster labels (col D) are re-ordered per sending cluster block,
# and all derived expression/specificity metrics (cols E-J, M-T) are refreshed
# with the new TPM-based figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "FAPs"  # D2
$ws.Cells.Item(2, 5).Value = 2  # E2
$ws.Cells.Item(2, 7).Value = 0.3541365  # G2
$ws.Cells.Item(2, 8).Value = 0.7082729999999999  # H2
$ws.Cells.Item(2, 9).Value = 0.06920610996416414  # I2
$ws.Cells.Item(2, 10).Value = 0.05608720193777861  # J2
$ws.Cells.Item(2, 13).Value = 220.2572175  # M2
$ws.Cells.Item(2, 14).Value = 440.514435  # N2
$ws.Cells.Item(2, 15).Value = 0.5129507012517006  # O2
$ws.Cells.Item(2, 16).Value = 0.4137724133030643  # P2
$ws.Cells.Item(2, 17).Value = 78.00112010518875  # Q2
$ws.Cells.Item(2, 18).Value = 312.004480420755  # R2
$ws.Cells.Item(2, 19).Value = 0.0354993226370203  # S2
$ws.Cells.Item(2, 20).Value = 0.02320733690121096  # T2

# Row 3
$ws.Cells.Item(3, 4).Value = "MuSCs"  # D3
$ws.Cells.Item(3, 5).Value = 2  # E3
$ws.Cells.Item(3, 7).Value = 0.3541365  # G3
$ws.Cells.Item(3, 8).Value = 0.7082729999999999  # H3
$ws.Cells.Item(3, 9).Value = 0.06920610996416414  # I3
$ws.Cells.Item(3, 10).Value = 0.05608720193777861  # J3
$ws.Cells.Item(3, 15).Value = 0.0001428715933923452  # O3
$ws.Cells.Item(3, 16).Value = 0.0001728713612618601  # P3
$ws.Cells.Item(3, 17).Value = 0.021725566002  # Q3
$ws.Cells.Item(3, 18).Value = 0.130353396012  # R3
$ws.Cells.Item(3, 19).Value = 0.000009887587203065986  # S3
$ws.Cells.Item(3, 20).Value = 0.000009695870948352626  # T3

# Row 4
$ws.Cells.Item(4, 4).Value = "ECs"  # D4
$ws.Cells.Item(4, 5).Value = 2  # E4
$ws.Cells.Item(4, 7).Value = 0.3541365  # G4
$ws.Cells.Item(4, 8).Value = 0.7082729999999999  # H4
$ws.Cells.Item(4, 9).Value = 0.06920610996416414  # I4
$ws.Cells.Item(4, 10).Value = 0.05608720193777861  # J4
$ws.Cells.Item(4, 13).Value = 71.69991033333334  # M4
$ws.Cells.Item(4, 14).Value = 215.099731  # N4
$ws.Cells.Item(4, 15).Value = 0.1669798597413381  # O4
$ws.Cells.Item(4, 16).Value = 0.202041812311349  # P4
$ws.Cells.Item(4, 17).Value = 25.3915552957605  # Q4
$ws.Cells.Item(4, 18).Value = 152.349331774563  # R4
$ws.Cells.Item(4, 19).Value = 0.01155602653505975  # S4
$ws.Cells.Item(4, 20).Value = 0.0113319599269814  # T4

# Row 5
$ws.Cells.Item(5, 4).Value = "Neutrophils"  # D5
$ws.Cells.Item(5, 5).Value = 2  # E5
$ws.Cells.Item(5, 7).Value = 0.3541365  # G5
$ws.Cells.Item(5, 8).Value = 0.7082729999999999  # H5
$ws.Cells.Item(5, 9).Value = 0.06920610996416414  # I5
$ws.Cells.Item(5, 10).Value = 0.05608720193777861  # J5
$ws.Cells.Item(5, 13).Value = 3.2906945  # M5
$ws.Cells.Item(5, 14).Value = 6.581389  # N5
$ws.Cells.Item(5, 15).Value = 0.007663603810758729  # O5
$ws.Cells.Item(5, 16).Value = 0.006181856922387211  # P5
$ws.Cells.Item(5, 17).Value = 1.16535503279925  # Q5
$ws.Cells.Item(5, 18).Value = 4.661420131197  # R5
$ws.Cells.Item(5, 19).Value = 0.000530368208049156  # S5
$ws.Cells.Item(5, 20).Value = 0.0003467230575563861  # T5

# Row 6
$ws.Cells.Item(6, 4).Value = "Inflammatory-Mac"  # D6
$ws.Cells.Item(6, 5).Value = 2  # E6
$ws.Cells.Item(6, 7).Value = 0.3541365  # G6
$ws.Cells.Item(6, 8).Value = 0.7082729999999999  # H6
$ws.Cells.Item(6, 9).Value = 0.06920610996416414  # I6
$ws.Cells.Item(6, 10).Value = 0.05608720193777861  # J6
$ws.Cells.Item(6, 13).Value = 27.607354  # M6
$ws.Cells.Item(6, 14).Value = 82.82206199999999  # N6
$ws.Cells.Item(6, 15).Value = 0.06429397299547716  # O6
$ws.Cells.Item(6, 16).Value = 0.07779423725008243  # P6
$ws.Cells.Item(6, 17).Value = 9.776771719820998  # Q6
$ws.Cells.Item(6, 18).Value = 58.66063031892599  # R6
$ws.Cells.Item(6, 19).Value = 0.004449535765157993  # S6
$ws.Cells.Item(6, 20).Value = 0.004363261094240832  # T6

# Row 7
$ws.Cells.Item(7, 5).Value = 2  # E7
$ws.Cells.Item(7, 7).Value = 0.3541365  # G7
$ws.Cells.Item(7, 8).Value = 0.7082729999999999  # H7
$ws.Cells.Item(7, 9).Value = 0.06920610996416414  # I7
$ws.Cells.Item(7, 10).Value = 0.05608720193777861  # J7
$ws.Cells.Item(7, 13).Value = 106.476041  # M7
$ws.Cells.Item(7, 14).Value = 319.428123  # N7
$ws.Cells.Item(7, 15).Value = 0.2479689906073331  # O7
$ws.Cells.Item(7, 16).Value = 0.3000368088518554  # P7
$ws.Cells.Item(7, 17).Value = 37.7070524935965  # Q7
$ws.Cells.Item(7, 18).Value = 226.242314961579  # R7
$ws.Cells.Item(7, 19).Value = 0.01716096923167388  # S7
$ws.Cells.Item(7, 20).Value = 0.01682822508684069  # T7

# Row 8
$ws.Cells.Item(8, 4).Value = "FAPs"  # D8
$ws.Cells.Item(8, 5).Value = 3  # E8
$ws.Cells.Item(8, 7).Value = 2.393812666666667  # G8
$ws.Cells.Item(8, 8).Value = 7.181438  # H8
$ws.Cells.Item(8, 9).Value = 0.4678039757069445  # I8
$ws.Cells.Item(8, 10).Value = 0.5686885753228443  # J8
$ws.Cells.Item(8, 13).Value = 220.2572175  # M8
$ws.Cells.Item(8, 14).Value = 440.514435  # N8
$ws.Cells.Item(8, 15).Value = 0.5129507012517006  # O8
$ws.Cells.Item(8, 16).Value = 0.4137724133030643  # P8
$ws.Cells.Item(8, 17).Value = 527.254517176255  # Q8
$ws.Cells.Item(8, 18).Value = 3163.52710305753  # R8
$ws.Cells.Item(8, 19).Value = 0.2399603773872107  # S8
$ws.Cells.Item(8, 20).Value = 0.2353076442292148  # T8

# Row 9
$ws.Cells.Item(9, 4).Value = "MuSCs"  # D9
$ws.Cells.Item(9, 5).Value = 3  # E9
$ws.Cells.Item(9, 7).Value = 2.393812666666667  # G9
$ws.Cells.Item(9, 8).Value = 7.181438  # H9
$ws.Cells.Item(9, 9).Value = 0.4678039757069445  # I9
$ws.Cells.Item(9, 10).Value = 0.5686885753228443  # J9
$ws.Cells.Item(9, 15).Value = 0.0001428715933923452  # O9
$ws.Cells.Item(9, 16).Value = 0.0001728713612618601  # P9
$ws.Cells.Item(9, 17).Value = 0.1468556194746667  # Q9
$ws.Cells.Item(9, 18).Value = 1.321700575272  # R9
$ws.Cells.Item(9, 19).Value = 0.00006683589940452509  # S9
$ws.Cells.Item(9, 20).Value = 0.00009830996815012795  # T9

# Row 10
$ws.Cells.Item(10, 4).Value = "ECs"  # D10
$ws.Cells.Item(10, 5).Value = 3  # E10
$ws.Cells.Item(10, 7).Value = 2.393812666666667  # G10
$ws.Cells.Item(10, 8).Value = 7.181438  # H10
$ws.Cells.Item(10, 9).Value = 0.4678039757069445  # I10
$ws.Cells.Item(10, 10).Value = 0.5686885753228443  # J10
$ws.Cells.Item(10, 13).Value = 71.69991033333334  # M10
$ws.Cells.Item(10, 14).Value = 215.099731  # N10
$ws.Cells.Item(10, 15).Value = 0.1669798597413381  # O10
$ws.Cells.Item(10, 16).Value = 0.202041812311349  # P10
$ws.Cells.Item(10, 17).Value = 171.6361535547976  # Q10
$ws.Cells.Item(10, 18).Value = 1544.725381993178  # R10
$ws.Cells.Item(10, 19).Value = 0.07811384224998591  # S10
$ws.Cells.Item(10, 20).Value = 0.1148988703989866  # T10

# Row 11
$ws.Cells.Item(11, 4).Value = "Neutrophils"  # D11
$ws.Cells.Item(11, 5).Value = 3  # E11
$ws.Cells.Item(11, 7).Value = 2.393812666666667  # G11
$ws.Cells.Item(11, 8).Value = 7.181438  # H11
$ws.Cells.Item(11, 9).Value = 0.4678039757069445  # I11
$ws.Cells.Item(11, 10).Value = 0.5686885753228443  # J11
$ws.Cells.Item(11, 13).Value = 3.2906945  # M11
$ws.Cells.Item(11, 14).Value = 6.581389  # N11
$ws.Cells.Item(11, 15).Value = 0.007663603810758729  # O11
$ws.Cells.Item(11, 16).Value = 0.006181856922387211  # P11
$ws.Cells.Item(11, 17).Value = 7.877306176230333  # Q11
$ws.Cells.Item(11, 18).Value = 47.263837057382  # R11
$ws.Cells.Item(11, 19).Value = 0.003585064330915824  # S11
$ws.Cells.Item(11, 20).Value = 0.003515551406042046  # T11

# Row 12
$ws.Cells.Item(12, 4).Value = "Inflammatory-Mac"  # D12
$ws.Cells.Item(12, 5).Value = 3  # E12
$ws.Cells.Item(12, 7).Value = 2.393812666666667  # G12
$ws.Cells.Item(12, 8).Value = 7.181438  # H12
$ws.Cells.Item(12, 9).Value = 0.4678039757069445  # I12
$ws.Cells.Item(12, 10).Value = 0.5686885753228443  # J12
$ws.Cells.Item(12, 13).Value = 27.607354  # M12
$ws.Cells.Item(12, 14).Value = 82.82206199999999  # N12
$ws.Cells.Item(12, 15).Value = 0.06429397299547716  # O12
$ws.Cells.Item(12, 16).Value = 0.07779423725008243  # P12
$ws.Cells.Item(12, 17).Value = 66.08683369835066  # Q12
$ws.Cells.Item(12, 18).Value = 594.781503285156  # R12
$ws.Cells.Item(12, 19).Value = 0.03007697618127914  # S12
$ws.Cells.Item(12, 20).Value = 0.04424069395007672  # T12

# Row 13
$ws.Cells.Item(13, 5).Value = 3  # E13
$ws.Cells.Item(13, 7).Value = 2.393812666666667  # G13
$ws.Cells.Item(13, 8).Value = 7.181438  # H13
$ws.Cells.Item(13, 9).Value = 0.4678039757069445  # I13
$ws.Cells.Item(13, 10).Value = 0.5686885753228443  # J13
$ws.Cells.Item(13, 13).Value = 106.476041  # M13
$ws.Cells.Item(13, 14).Value = 319.428123  # N13
$ws.Cells.Item(13, 15).Value = 0.2479689906073331  # O13
$ws.Cells.Item(13, 16).Value = 0.3000368088518554  # P13
$ws.Cells.Item(13, 17).Value = 254.8836956423194  # Q13
$ws.Cells.Item(13, 18).Value = 2293.953260780874  # R13
$ws.Cells.Item(13, 19).Value = 0.1160008796581484  # S13
$ws.Cells.Item(13, 20).Value = 0.1706275053703742  # T13

# Row 14
$ws.Cells.Item(14, 4).Value = "FAPs"  # D14
$ws.Cells.Item(14, 6).Value = 1  # F14
$ws.Cells.Item(14, 7).Value = 2.3691785  # G14
$ws.Cells.Item(14, 8).Value = 4.738357  # H14
$ws.Cells.Item(14, 9).Value = 0.4629899143288914  # I14
$ws.Cells.Item(14, 10).Value = 0.3752242227393771  # J14
$ws.Cells.Item(14, 13).Value = 220.2572175  # M14
$ws.Cells.Item(14, 14).Value = 440.514435  # N14
$ws.Cells.Item(14, 15).Value = 0.5129507012517006  # O14
$ws.Cells.Item(14, 16).Value = 0.4137724133030643  # P14
$ws.Cells.Item(14, 17).Value = 521.8286641708237  # Q14
$ws.Cells.Item(14, 18).Value = 2087.314656683295  # R14
$ws.Cells.Item(14, 19).Value = 0.2374910012274696  # S14
$ws.Cells.Item(14, 20).Value = 0.1552574321726386  # T14

# Row 15
$ws.Cells.Item(15, 4).Value = "MuSCs"  # D15
$ws.Cells.Item(15, 6).Value = 1  # F15
$ws.Cells.Item(15, 7).Value = 2.3691785  # G15
$ws.Cells.Item(15, 8).Value = 4.738357  # H15
$ws.Cells.Item(15, 9).Value = 0.4629899143288914  # I15
$ws.Cells.Item(15, 10).Value = 0.3752242227393771  # J15
$ws.Cells.Item(15, 15).Value = 0.0001428715933923452  # O15
$ws.Cells.Item(15, 16).Value = 0.0001728713612618601  # P15
$ws.Cells.Item(15, 17).Value = 0.145344362618  # Q15
$ws.Cells.Item(15, 18).Value = 0.872066175708  # R15
$ws.Cells.Item(15, 19).Value = 0.00006614810678475408  # S15
$ws.Cells.Item(15, 20).Value = 0.00006486552216337952  # T15

# Row 16
$ws.Cells.Item(16, 4).Value = "ECs"  # D16
$ws.Cells.Item(16, 6).Value = 1  # F16
$ws.Cells.Item(16, 7).Value = 2.3691785  # G16
$ws.Cells.Item(16, 8).Value = 4.738357  # H16
$ws.Cells.Item(16, 9).Value = 0.4629899143288914  # I16
$ws.Cells.Item(16, 10).Value = 0.3752242227393771  # J16
$ws.Cells.Item(16, 13).Value = 71.69991033333334  # M16
$ws.Cells.Item(16, 14).Value = 215.099731  # N16
$ws.Cells.Item(16, 15).Value = 0.1669798597413381  # O16
$ws.Cells.Item(16, 16).Value = 0.202041812311349  # P16
$ws.Cells.Item(16, 17).Value = 169.8698860136612  # Q16
$ws.Cells.Item(16, 18).Value = 1019.219316081967  # R16
$ws.Cells.Item(16, 19).Value = 0.0773099909562924  # S16
$ws.Cells.Item(16, 20).Value = 0.07581098198538105  # T16

# Row 17
$ws.Cells.Item(17, 4).Value = "Neutrophils"  # D17
$ws.Cells.Item(17, 6).Value = 1  # F17
$ws.Cells.Item(17, 7).Value = 2.3691785  # G17
$ws.Cells.Item(17, 8).Value = 4.738357  # H17
$ws.Cells.Item(17, 9).Value = 0.4629899143288914  # I17
$ws.Cells.Item(17, 10).Value = 0.3752242227393771  # J17
$ws.Cells.Item(17, 13).Value = 3.2906945  # M17
$ws.Cells.Item(17, 14).Value = 6.581389  # N17
$ws.Cells.Item(17, 15).Value = 0.007663603810758729  # O17
$ws.Cells.Item(17, 16).Value = 0.006181856922387211  # P17
$ws.Cells.Item(17, 17).Value = 7.796242659468249  # Q17
$ws.Cells.Item(17, 18).Value = 31.184970637873  # R17
$ws.Cells.Item(17, 19).Value = 0.003548171271793749  # S17
$ws.Cells.Item(17, 20).Value = 0.002319582458788779  # T17

# Row 18
$ws.Cells.Item(18, 4).Value = "Inflammatory-Mac"  # D18
$ws.Cells.Item(18, 6).Value = 1  # F18
$ws.Cells.Item(18, 7).Value = 2.3691785  # G18
$ws.Cells.Item(18, 8).Value = 4.738357  # H18
$ws.Cells.Item(18, 9).Value = 0.4629899143288914  # I18
$ws.Cells.Item(18, 10).Value = 0.3752242227393771  # J18
$ws.Cells.Item(18, 13).Value = 27.607354  # M18
$ws.Cells.Item(18, 14).Value = 82.82206199999999  # N18
$ws.Cells.Item(18, 15).Value = 0.06429397299547716  # O18
$ws.Cells.Item(18, 16).Value = 0.07779423725008243  # P18
$ws.Cells.Item(18, 17).Value = 65.40674953868898  # Q18
$ws.Cells.Item(18, 18).Value = 392.4404972321339  # R18
$ws.Cells.Item(18, 19).Value = 0.02976746104904003  # S18
$ws.Cells.Item(18, 20).Value = 0.02919028220576488  # T18

# Row 19
$ws.Cells.Item(19, 6).Value = 1  # F19
$ws.Cells.Item(19, 7).Value = 2.3691785  # G19
$ws.Cells.Item(19, 8).Value = 4.738357  # H19
$ws.Cells.Item(19, 9).Value = 0.4629899143288914  # I19
$ws.Cells.Item(19, 10).Value = 0.3752242227393771  # J19
$ws.Cells.Item(19, 13).Value = 106.476041  # M19
$ws.Cells.Item(19, 14).Value = 319.428123  # N19
$ws.Cells.Item(19, 15).Value = 0.2479689906073331  # O19
$ws.Cells.Item(19, 16).Value = 0.3000368088518554  # P19
$ws.Cells.Item(19, 17).Value = 252.2607471023185  # Q19
$ws.Cells.Item(19, 18).Value = 1513.564482613911  # R19
$ws.Cells.Item(19, 19).Value = 0.1148071417175108  # S19
$ws.Cells.Item(19, 20).Value = 0.1125810783946405  # T19
